# Update Azure AD group object IDs on the "General" sheet for the new tenant.
# General!B10 -> dcAdminsAadGroupObjectId   (used as dcSubnetVmContAadGroupId)
# General!B11 -> serverTeamAadGroupObjectId (used as hubSubnetVmContAadGroupId)
# General!B12 -> appTeamsAadGroupObjectId   (used as spokeVmContAadGroupId)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

$ws.Range("B10").Value = "83578c91-9919-4bd8-bee8-2649f6eb7c13"
$ws.Range("B11").Value = "57f2ff92-300b-4075-a7ab-2030b46ebe2f"
$ws.Range("B12").Value = "002984bd-b5ce-445d-8138-d19b514550c7"

# Reflect the user's selection after editing those three cells.
$ws.Range("B10:B12").Select()
